$wb = $excel.ActiveWorkbook

# --- Sheet "Home win": add a new prediction row (row 3) ---
$wsHome = $wb.Worksheets.Item("Home win")
$wsHome.Range("A3").Value = "04-02-2025 20:00"
$wsHome.Range("B3").Value = "ENGLAND"
$wsHome.Range("C3").Value = "PREMIER LEAGUE CUP"
$wsHome.Range("D3").Value = "Ipswich Town U21 - Watford U21"
$wsHome.Range("E3").Value = 73.3
$wsHome.Range("F3").Value = 1.73

# --- Sheet "Btts": remove stale predictions (old rows 3-5), keep the rest ---
$wsBtts = $wb.Worksheets.Item("Btts")
$wsBtts.Rows("3:5").Delete()

# --- Sheet "Over_Under": remove a stale prediction (old row 3) ---
$wsOverUnder = $wb.Worksheets.Item("Over_Under")
$wsOverUnder.Rows("3:3").Delete()
